$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. A leading apostrophe forces Excel to
# treat the value as literal text (these price/volume columns are plain
# strings in the source data, some of which look like numbers or dates).
$updates = [ordered]@{
    "D2" = "26.291.80"
    "E2" = "  +1.11%  "
    "D3" = "1.677.65"
    "E3" = "  +0.78%  "
    "D4" = "1.009"
    "E4" = "  +0.39%  "
    "D5" = "217.65"
    "E5" = "  +0.53%  "
    "D6" = "0.5254"
    "E6" = "  +3.20%  "
    "D7" = "1.009"
    "E7" = "  +0.32%  "
    "D8" = "0.2685"
    "E8" = "  +2.09%  "
    "D9" = "0.06470"
    "E9" = "  +1.34%  "
    "D10" = "21.87"
    "E10" = "  +0.42%  "
    "D11" = "0.07509"
    "E11" = "  +1.45%  "
    "D12" = "1.704.04"
    "E12" = "  +2.34%  "
    "D13" = "4.512"
    "E13" = "  +0.34%  "
    "D14" = "0.5766"
    "E14" = "  -0.78%  "
    "D15" = "0.000008469"
    "E15" = "  -0.21%  "
    "D16" = "64.63"
    "E16" = "  +0.74%  "
    "D17" = "26.318.56"
    "E17" = "  +1.00%  "
    "D18" = "4.918"
    "E18" = "  +0.32%  "
    "D19" = "1.009"
    "E19" = "  +0.35%  "
    "D20" = "10.86"
    "E20" = "  +1.85%  "
    "D21" = "189.79"
    "E21" = "  +0.48%  "
    "D22" = "6.181"
    "E22" = "  -0.21%  "
    "E23" = "  +0.37%  "
    "D24" = "144.87"
    "E24" = "  -0.28%  "
    "B25" = "Cosmos"
    "C25" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D25" = "7.795"
    "E25" = "  +2.92%  "
    "B26" = "Stellar"
    "C26" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D26" = "0.1261"
    "E26" = "  +6.19%  "
    "D27" = "15.74"
    "E27" = "  +0.88%  "
    "D28" = "0.06414"
    "E28" = "  -3.87%  "
    "D29" = "1.363"
    "E29" = "  +4.75%  "
    "E30" = "  +0.38%  "
    "D31" = "3.580"
    "E31" = "  +1.67%  "
    "D32" = "3.579"
    "E32" = "  +2.40%  "
    "D33" = "1.654"
    "E33" = "  +1.88%  "
    "D34" = "1.026"
    "E34" = "  +0.94%  "
    "D35" = "0.6179"
    "E35" = "  +1.98%  "
    "D36" = "2.407"
    "D37" = "2.742"
    "E37" = "  +2.23%  "
    "D38" = "6.283"
    "E38" = "  +1.16%  "
    "D39" = "1.116.82"
    "E39" = "  +3.85%  "
    "D40" = "0.01621"
    "E40" = "  +0.75%  "
    "D41" = "0.8710"
    "E41" = "  +1.44%  "
    "D42" = "1.016"
    "E42" = "  +0.74%  "
    "E43" = "  +0.15%  "
    "D44" = "1.827.70"
    "E44" = "  +0.89%  "
    "B45" = "BabyDogeCoin"
    "C45" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D45" = "0.00000000110"
    "E45" = "  -2.52%  "
    "B46" = "Aave"
    "C46" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D46" = "56.90"
    "E46" = "  +1.25%  "
    "B47" = "Frax"
    "C47" = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
    "D47" = "1.008"
    "E47" = "  +0.28%  "
    "B48" = "EnergySwap"
    "C48" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D48" = "8.167"
    "E48" = "  +1.86%  "
    "E49" = "  +1.11%  "
    "D50" = "0.4299"
    "E50" = "  +0.19%  "
    "D51" = "6.052"
    "E51" = "  +1.89%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    # Setting the value via the quote-prefix above can stamp a transient
    # "text" number format onto the cell; restore the Normal style so the
    # cell keeps its original (unstyled) appearance.
    $cell.Style = "Normal"
}
